# Settings back to qfin22 presentation - recalibration
$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.0001993655275788096
$ws.Range("B3").Value = 0.7932514652920688
$ws.Range("B4").Value = 0.001501574309577561
$ws.Range("B5").Value = "[1.0, 0.3428654924548603, 0.334339821203891, 0.3107505615876993, 0.301019256423924, 0.4983087885708196, 0.32541569437551643, 0.330259293569174, 0.31308677919220035, 0.2914262482884507, 0.28683857327868606, 0.27634481670143557, 0.28952962931217024, 0.3062349194322275, 0.32067672151483423, 0.2949252840086501, 0.2714451057352363, 0.285960540053174, 0.2817363883217259, 0.27848678560222834]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.009526733307493023
$ws.Range("B5").Value = 0.7653730596278702
$ws.Range("B6").Value = 0.0009540231914201716
$ws.Range("B7").Value = 0.009929380664100141
$ws.Range("B8").Value = 0.758539439045549
$ws.Range("B9").Value = 0.001034520974348429
$ws.Range("B10").Value = "[1.0, 0.36379040018019687, 0.3369116524599222, 0.33423415324670563, 0.30221423163808236, 0.37065897584399504, 0.31187706628447076, 0.3076954202528062, 0.29278294306587144, 0.2686672025436916, 0.2943322687329038, 0.26124095902847033, 0.2628770400642022, 0.25636746519986403, 0.28275826869499066, 0.30686152475512857, 0.2605032607586898, 0.2626711409502044, 0.27283393942022244, 0.2561667741484799]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0003275724123776003
$ws.Range("B3").Value = 0.00003349873430469541
$ws.Range("B4").Value = 0.1999997027857356
$ws.Range("B5").Value = 0.7800007457931163
$ws.Range("B6").Value = "[1.0, -0.08168296419615892, -0.049041386496671086, -0.09597174397217319, -0.09298087705772949, 0.29885021646636845, -0.031686367013017344, 0.06344537085548103, -0.0022526950692061285, -0.010440818243234888, -0.028540968437931837, -0.025864123779673356, 0.003661736680930923, 0.03432216793854821, 0.07399842013413903, 0.023886402257475306, 0.004687550904959184, -0.00912806927556006, -0.01724756621995928, 0.017067800127909736]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0005209117223476958
$ws.Range("B3").Value = 0.00003349715710288374
$ws.Range("B4").Value = 0.1000000684366096
$ws.Range("B5").Value = 0.8299999368862377
$ws.Range("B6").Value = "[1.0, -0.04802030782139957, -0.023900929377378574, -0.0743374040175018, -0.0705543274961626, 0.3097107309070136, -0.017452237246948558, 0.0733423485753081, 0.008217110779852573, -0.004166913578727237, -0.020504306001846613, -0.02155455602500132, 0.00952647179665617, 0.03935998615046954, 0.07636647426316716, 0.030319645173745163, 0.006183373235456205, -0.004567785579492703, -0.012508178077142149, 0.02353353719599398]"
$ws.Range("B7").Value = 0.09999991635525493

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.00007872550446405591
$ws.Range("B3").Value = 0.00003018114400777296
$ws.Range("B4").Value = 0.2000153345572102
$ws.Range("B5").Value = 0.7750091334972046
$ws.Range("B6").Value = "[1.0, -0.06478753956565701, -0.024912513191969073, -0.06269228363032286, -0.04605218636580197, 0.21817205258916988, -0.02478730161821368, 0.039702398133601734, -0.005871947513234957, -0.0013599584585808253, -0.013409119946626113, -0.015087497781139369, 0.0009088577674852645, 0.020455933716936336, 0.06602904454261269, 0.025897981234409805, 0.0020496803740213737, 0.005066724737274858, 0.011462502315722928, 0.01699976373146314]"
$ws.Range("B7").Value = 0.0100022938631411
$ws.Range("B9").Value = 0.7713498724234
